$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value is a plain numeric-looking string (e.g. "100.41")
# must be forced to remain text, matching the source workbook where every
# data cell is stored as a literal string. A leading apostrophe is the
# standard Excel "treat as text" input prefix and keeps all of these cells
# sharing a single extra (quote-prefixed) style, instead of drifting in value.
$textCells = @("D6", "D7", "D9", "D10", "D11", "D15", "D18", "D21", "D22", "D23", "D25", "D26", "D27", "D29", "D30", "D32", "D33", "D34", "D40", "D41", "D42", "D43", "D45", "D47", "D49", "D50", "D51")

function Set-CellText($range, $value) {
    if ($textCells -contains $range) {
        $ws.Range($range).Value = "'" + $value
    } else {
        $ws.Range($range).Value = $value
    }
}

Set-CellText "D2" "42.508.57"
Set-CellText "E2" "  -2.35%  "
Set-CellText "D3" "2.353.87"
Set-CellText "E3" "  +0.20%  "
Set-CellText "E4" "  -0.31%  "
Set-CellText "E5" "  +3.84%  "
Set-CellText "D6" "100.41"
Set-CellText "E6" "  -8.52%  "
Set-CellText "D7" "0.637"
Set-CellText "E7" "  -1.33%  "
Set-CellText "E8" "  -0.06%  "
Set-CellText "D9" "0.622"
Set-CellText "E9" "  -2.17%  "
Set-CellText "D10" "39.84"
Set-CellText "E10" "  -7.89%  "
Set-CellText "D11" "0.0921"
Set-CellText "E11" "  -2.02%  "
Set-CellText "E12" "  -4.89%  "
Set-CellText "E13" "  -3.49%  "
Set-CellText "E14" "  +0.02%  "
Set-CellText "D15" "16.48"
Set-CellText "E15" "  +0.70%  "
Set-CellText "D16" "2.709.35"
Set-CellText "E16" "  -0.04%  "
Set-CellText "D17" "2.354.31"
Set-CellText "E17" "  -3.06%  "
Set-CellText "D18" "8.05"
Set-CellText "E18" "  +11.03%  "
Set-CellText "D19" "42.526.74"
Set-CellText "E19" "  -2.30%  "
Set-CellText "E20" "  -2.09%  "
Set-CellText "D21" "76.27"
Set-CellText "E21" "  +0.73%  "
Set-CellText "D22" "3.70"
Set-CellText "E22" "  +7.41%  "
Set-CellText "D23" "265.36"
Set-CellText "E23" "  +3.36%  "
Set-CellText "E24" "  -10.22%  "
Set-CellText "D25" "10.01"
Set-CellText "E25" "  +9.53%  "
Set-CellText "D26" "1.00"
Set-CellText "E26" "  +0.08%  "
Set-CellText "D27" "11.43"
Set-CellText "E27" "  -5.54%  "
Set-CellText "E28" "  +2.30%  "
Set-CellText "D29" "2.21"
Set-CellText "E29" "  -1.97%  "
Set-CellText "D30" "175.10"
Set-CellText "E30" "  +0.72%  "
Set-CellText "E31" "  -2.51%  "
Set-CellText "D32" "0.0897"
Set-CellText "E32" "  -3.78%  "
Set-CellText "D33" "35.07"
Set-CellText "E33" "  -10.60%  "
Set-CellText "D34" "6.02"
Set-CellText "E34" "  -0.18%  "
Set-CellText "E35" "  -0.39%  "
Set-CellText "E36" "  -8.55%  "
Set-CellText "E37" "  -5.22%  "
Set-CellText "E38" "  +8.25%  "
Set-CellText "E39" "  +0.96%  "
Set-CellText "D40" "3.78"
Set-CellText "E40" "  -9.19%  "
Set-CellText "D41" "1.50"
Set-CellText "E41" "  +1.42%  "
Set-CellText "D42" "0.234"
Set-CellText "E42" "  -0.10%  "
Set-CellText "D43" "69.86"
Set-CellText "E43" "  -3.58%  "
Set-CellText "E44" "  -0.28%  "
Set-CellText "D45" "118.99"
Set-CellText "E45" "  +6.99%  "
Set-CellText "E46" "  +20.73%  "
Set-CellText "D47" "11.86"
Set-CellText "E47" "  -7.39%  "
Set-CellText "E48" "  -2.66%  "
Set-CellText "D49" "9.19"
Set-CellText "E49" "  -1.12%  "
Set-CellText "B50" "ordi"
Set-CellText "C50" "https://coinranking.com/coin/j7-7vPrOi+ordi-ordi"
Set-CellText "D50" "72.80"
Set-CellText "E50" "  +3.64%  "
Set-CellText "B51" "TrustWalletToken"
Set-CellText "C51" "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
Set-CellText "D51" "1.26"
Set-CellText "E51" "  -4.30%  "
